$wb = $excel.ActiveWorkbook

# --- 1) Rename the "Requested quantity" header to the new metric names ----
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2) Add the new "PO Forecast" sheet after the last existing sheet -----
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# --- 3) Header row: ds / PO_Forecast / yhat_lower / yhat_upper ------------
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Re-use the bold/bordered header style already used on the other sheets
$wsWeekly.Range("A1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# --- 4) Data rows 2..83 -----------------------------------------------------
$arr = New-Object 'object[,]' 82,4
$arr[0,0]=44934.99999999999; $arr[0,1]=142; $arr[0,2]=-245.1818025825996; $arr[0,3]=564.9061943936537
$arr[1,0]=44941.99999999999; $arr[1,1]=146; $arr[1,2]=-239.3052679417948; $arr[1,3]=568.6634497758016
$arr[2,0]=44948.99999999999; $arr[2,1]=150; $arr[2,2]=-271.842556229648; $arr[2,3]=543.972017898505
$arr[3,0]=44955.99999999999; $arr[3,1]=154; $arr[3,2]=-301.5136175832059; $arr[3,3]=566.9742394448638
$arr[4,0]=44969.99999999999; $arr[4,1]=162; $arr[4,2]=-211.7361246213262; $arr[4,3]=602.3283556378829
$arr[5,0]=44976.99999999999; $arr[5,1]=166; $arr[5,2]=-225.9648177095376; $arr[5,3]=556.2771248927088
$arr[6,0]=44983.99999999999; $arr[6,1]=170; $arr[6,2]=-223.4992150723068; $arr[6,3]=594.2832528839565
$arr[7,0]=44990.99999999999; $arr[7,1]=174; $arr[7,2]=-262.9882303291892; $arr[7,3]=563.8930749239114
$arr[8,0]=44997.99999999999; $arr[8,1]=178; $arr[8,2]=-230.0735876822071; $arr[8,3]=590.9877013576095
$arr[9,0]=45004.99999999999; $arr[9,1]=182; $arr[9,2]=-241.4262126136746; $arr[9,3]=585.6895053579354
$arr[10,0]=45011.99999999999; $arr[10,1]=186; $arr[10,2]=-226.3514344196527; $arr[10,3]=561.305448153925
$arr[11,0]=45025.99999999999; $arr[11,1]=194; $arr[11,2]=-208.8214119561563; $arr[11,3]=617.429949923956
$arr[12,0]=45032.99999999999; $arr[12,1]=198; $arr[12,2]=-221.7401894242824; $arr[12,3]=626.2200794316943
$arr[13,0]=45039.99999999999; $arr[13,1]=202; $arr[13,2]=-198.9860238080499; $arr[13,3]=613.1126843538657
$arr[14,0]=45046.99999999999; $arr[14,1]=206; $arr[14,2]=-225.1998793714144; $arr[14,3]=609.8786334283412
$arr[15,0]=45053.99999999999; $arr[15,1]=210; $arr[15,2]=-187.6157208001995; $arr[15,3]=605.6472569536071
$arr[16,0]=45060.99999999999; $arr[16,1]=214; $arr[16,2]=-202.1015754793416; $arr[16,3]=641.2952415645882
$arr[17,0]=45067.99999999999; $arr[17,1]=218; $arr[17,2]=-184.8751283053966; $arr[17,3]=625.0171528951728
$arr[18,0]=45074.99999999999; $arr[18,1]=222; $arr[18,2]=-175.2885415033917; $arr[18,3]=600.3401655277472
$arr[19,0]=45081.99999999999; $arr[19,1]=226; $arr[19,2]=-181.0019956140603; $arr[19,3]=593.5741296815593
$arr[20,0]=45088.99999999999; $arr[20,1]=230; $arr[20,2]=-177.5971397170799; $arr[20,3]=665.7215136152141
$arr[21,0]=45095.99999999999; $arr[21,1]=234; $arr[21,2]=-176.8801874115024; $arr[21,3]=609.3072894348401
$arr[22,0]=45102.99999999999; $arr[22,1]=238; $arr[22,2]=-162.6696166860614; $arr[22,3]=613.8410667646483
$arr[23,0]=45109.99999999999; $arr[23,1]=242; $arr[23,2]=-159.9412919393766; $arr[23,3]=639.4847656964693
$arr[24,0]=45116.99999999999; $arr[24,1]=246; $arr[24,2]=-168.8237570444087; $arr[24,3]=660.9730661993476
$arr[25,0]=45123.99999999999; $arr[25,1]=250; $arr[25,2]=-152.2384179664733; $arr[25,3]=653.8799854889711
$arr[26,0]=45130.99999999999; $arr[26,1]=254; $arr[26,2]=-163.6662159609827; $arr[26,3]=661.9947185849103
$arr[27,0]=45137.99999999999; $arr[27,1]=258; $arr[27,2]=-122.5280250662382; $arr[27,3]=687.4589124301591
$arr[28,0]=45151.99999999999; $arr[28,1]=266; $arr[28,2]=-106.7991275521457; $arr[28,3]=670.8192577725495
$arr[29,0]=45172.99999999999; $arr[29,1]=278; $arr[29,2]=-127.4712768521209; $arr[29,3]=701.6903035224941
$arr[30,0]=45179.99999999999; $arr[30,1]=282; $arr[30,2]=-119.9422441310587; $arr[30,3]=670.0909973596792
$arr[31,0]=45186.99999999999; $arr[31,1]=286; $arr[31,2]=-139.8802088526832; $arr[31,3]=705.2777075656877
$arr[32,0]=45193.99999999999; $arr[32,1]=290; $arr[32,2]=-107.0811184137038; $arr[32,3]=702.1308660180467
$arr[33,0]=45200.99999999999; $arr[33,1]=294; $arr[33,2]=-104.7570071449671; $arr[33,3]=699.0363774453651
$arr[34,0]=45228.99999999999; $arr[34,1]=310; $arr[34,2]=-116.587441950214; $arr[34,3]=720.9982779720104
$arr[35,0]=45235.99999999999; $arr[35,1]=314; $arr[35,2]=-94.54126658181423; $arr[35,3]=697.3339362757225
$arr[36,0]=45242.99999999999; $arr[36,1]=318; $arr[36,2]=-76.92039270011524; $arr[36,3]=730.8706042820584
$arr[37,0]=45249.99999999999; $arr[37,1]=322; $arr[37,2]=-92.74536743034535; $arr[37,3]=721.1825731593442
$arr[38,0]=45256.99999999999; $arr[38,1]=326; $arr[38,2]=-72.7096256565593; $arr[38,3]=729.8942189430238
$arr[39,0]=45263.99999999999; $arr[39,1]=330; $arr[39,2]=-73.23020700059332; $arr[39,3]=715.9866534533779
$arr[40,0]=45270.99999999999; $arr[40,1]=334; $arr[40,2]=-67.38775122420883; $arr[40,3]=725.8097380739562
$arr[41,0]=45277.99999999999; $arr[41,1]=338; $arr[41,2]=-55.401821695193; $arr[41,3]=752.1809040232724
$arr[42,0]=45298.99999999999; $arr[42,1]=350; $arr[42,2]=-67.57990777526831; $arr[42,3]=726.1329061812247
$arr[43,0]=45305.99999999999; $arr[43,1]=354; $arr[43,2]=-72.51268336064997; $arr[43,3]=759.5320403826239
$arr[44,0]=45319.99999999999; $arr[44,1]=362; $arr[44,2]=-32.99187096957441; $arr[44,3]=790.2577250587934
$arr[45,0]=45326.99999999999; $arr[45,1]=366; $arr[45,2]=-8.604138588427398; $arr[45,3]=767.6555216911248
$arr[46,0]=45333.99999999999; $arr[46,1]=370; $arr[46,2]=-45.10357526802417; $arr[46,3]=756.924358554459
$arr[47,0]=45340.99999999999; $arr[47,1]=374; $arr[47,2]=-16.81269946201979; $arr[47,3]=765.3945073117931
$arr[48,0]=45347.99999999999; $arr[48,1]=378; $arr[48,2]=-27.31154178187116; $arr[48,3]=806.3174200116923
$arr[49,0]=45354.99999999999; $arr[49,1]=382; $arr[49,2]=-32.40744981978147; $arr[49,3]=792.5333682320403
$arr[50,0]=45361.99999999999; $arr[50,1]=386; $arr[50,2]=15.86083829213782; $arr[50,3]=819.6863110869571
$arr[51,0]=45368.99999999999; $arr[51,1]=390; $arr[51,2]=-12.44071394731417; $arr[51,3]=795.1352764237735
$arr[52,0]=45375.99999999999; $arr[52,1]=394; $arr[52,2]=-6.231727546854906; $arr[52,3]=800.4883303225313
$arr[53,0]=45382.99999999999; $arr[53,1]=398; $arr[53,2]=-14.38962970346762; $arr[53,3]=798.8138118563746
$arr[54,0]=45389.99999999999; $arr[54,1]=402; $arr[54,2]=-21.41209515872223; $arr[54,3]=763.7602381631131
$arr[55,0]=45396.99999999999; $arr[55,1]=406; $arr[55,2]=24.35991846856425; $arr[55,3]=811.0312466769732
$arr[56,0]=45410.99999999999; $arr[56,1]=414; $arr[56,2]=2.34705065576923; $arr[56,3]=823.6411585717143
$arr[57,0]=45417.99999999999; $arr[57,1]=418; $arr[57,2]=-10.13808851515972; $arr[57,3]=837.6398504930316
$arr[58,0]=45424.99999999999; $arr[58,1]=422; $arr[58,2]=14.58054145491072; $arr[58,3]=804.0005467191955
$arr[59,0]=45431.99999999999; $arr[59,1]=426; $arr[59,2]=33.83092017558202; $arr[59,3]=817.0324935201406
$arr[60,0]=45438.99999999999; $arr[60,1]=430; $arr[60,2]=32.68956474352508; $arr[60,3]=809.5795588962532
$arr[61,0]=45445.99999999999; $arr[61,1]=434; $arr[61,2]=0.4805319514152343; $arr[61,3]=812.098321061637
$arr[62,0]=45452.99999999999; $arr[62,1]=438; $arr[62,2]=60.63008279796053; $arr[62,3]=835.0905603827168
$arr[63,0]=45459.99999999999; $arr[63,1]=442; $arr[63,2]=34.55233906634367; $arr[63,3]=859.227355529304
$arr[64,0]=45466.99999999999; $arr[64,1]=446; $arr[64,2]=56.1703306138925; $arr[64,3]=870.2846771724204
$arr[65,0]=45494.99999999999; $arr[65,1]=462; $arr[65,2]=59.29399125637091; $arr[65,3]=874.681910790685
$arr[66,0]=45501.99999999999; $arr[66,1]=466; $arr[66,2]=61.57185963552685; $arr[66,3]=897.4118036650648
$arr[67,0]=45508.99999999999; $arr[67,1]=470; $arr[67,2]=65.23410271628427; $arr[67,3]=876.8155998183129
$arr[68,0]=45543.99999999999; $arr[68,1]=490; $arr[68,2]=65.94368405873578; $arr[68,3]=910.8730523217528
$arr[69,0]=45550.99999999999; $arr[69,1]=494; $arr[69,2]=54.06030864762864; $arr[69,3]=892.7855432928587
$arr[70,0]=45557.99999999999; $arr[70,1]=498; $arr[70,2]=92.05979612700824; $arr[70,3]=920.2688625723267
$arr[71,0]=45578.99999999999; $arr[71,1]=510; $arr[71,2]=91.87931565949205; $arr[71,3]=935.8710229693809
$arr[72,0]=45599.99999999999; $arr[72,1]=522; $arr[72,2]=128.7603567256844; $arr[72,3]=959.1316360222828
$arr[73,0]=45620.99999999999; $arr[73,1]=534; $arr[73,2]=135.6040548908716; $arr[73,3]=938.3243702723009
$arr[74,0]=45627.99999999999; $arr[74,1]=539; $arr[74,2]=131.4253154286945; $arr[74,3]=958.057636037012
$arr[75,0]=45634.99999999999; $arr[75,1]=543; $arr[75,2]=163.9279740841779; $arr[75,3]=982.2108666700436
$arr[76,0]=45641.99999999999; $arr[76,1]=547; $arr[76,2]=135.6248489029222; $arr[76,3]=946.4398993285247
$arr[77,0]=45648.99999999999; $arr[77,1]=551; $arr[77,2]=156.5516586298905; $arr[77,3]=937.5880184717632
$arr[78,0]=45655.99999999999; $arr[78,1]=555; $arr[78,2]=133.9904199085311; $arr[78,3]=957.8164479011909
$arr[79,0]=45662.99999999999; $arr[79,1]=559; $arr[79,2]=156.3733459522105; $arr[79,3]=996.4416003408297
$arr[80,0]=45669.99999999999; $arr[80,1]=563; $arr[80,2]=157.7245326052636; $arr[80,3]=971.4934940466129
$arr[81,0]=45676.99999999999; $arr[81,1]=567; $arr[81,2]=187.3964844471547; $arr[81,3]=970.9067911032015

$wsForecast.Range("A2:D83").Value = $arr

# Re-use the date-style formatting (column A) from the other sheets
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A83").PasteSpecial(-4122)

$wsForecast.Range("A1").Select()

Write-Host "PO Forecast sheet added with $($wsForecast.Range('A2:D83').Rows.Count) data rows"
